$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the existing sheets.
# ---------------------------------------------------------------------------
$wsName = $wb.Worksheets.Item(1)
$wsName.Name = "NameEntities"

$wsDesc = $wb.Worksheets.Item(2)
$wsDesc.Name = "DescriptionEntities"

# ---------------------------------------------------------------------------
# 2. Add the two new quest localization sheets at the end, in order.
# ---------------------------------------------------------------------------
$wsQuestName = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsQuestName.Name = "QuestNameEntities"

$wsQuestDesc = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsQuestDesc.Name = "QuestDescriptionEntities"

# ---------------------------------------------------------------------------
# 3. Populate QuestDescriptionEntities header + first english rows (this is
#    the order the original author typed things in, which drives the shared
#    string table order).
# ---------------------------------------------------------------------------
$wsQuestDesc.Range("A1").Value = "id"
$wsQuestDesc.Range("B1").Value = "en01"
$wsQuestDesc.Range("C1").Value = "en02"
$wsQuestDesc.Range("D1").Value = "en03"
$wsQuestDesc.Range("E1").Value = "jp01"
$wsQuestDesc.Range("F1").Value = "jp02"
$wsQuestDesc.Range("G1").Value = "jp03"

$wsQuestDesc.Range("A2").Value = 20000
$wsQuestDesc.Range("B2").Value = "Call my father?"

$wsQuestDesc.Range("A3").Value = 25000
$wsQuestDesc.Range("B3").Value = "Scavenge for coins."
$wsQuestDesc.Range("C3").Value = "Purchase Amulet from the vending machine"

# ---------------------------------------------------------------------------
# 4. Populate QuestNameEntities completely (english column fully, then the
#    japanese column, matching the original authoring order).
# ---------------------------------------------------------------------------
$wsQuestName.Range("A1").Value = "id"
$wsQuestName.Range("B1").Value = "en"
$wsQuestName.Range("C1").Value = "jp"

$wsQuestName.Range("A2").Value = 20000
$wsQuestName.Range("A3").Value = 25000

$wsQuestName.Range("B2").Value = "An unknown presence.."
$wsQuestName.Range("B3").Value = "Economic hardship"

$wsQuestName.Range("C2").Value = "未知の存在"
$wsQuestName.Range("C3").Value = "びんぼう"

# ---------------------------------------------------------------------------
# 5. Back to QuestDescriptionEntities for the japanese translations.
# ---------------------------------------------------------------------------
$wsQuestDesc.Range("E2").Value = "パパを連絡"

$wsQuestDesc.Range("E3").Value = "小銭を探せ"
$wsQuestDesc.Range("F3").Value = "自動販売機でアミュレットを購入する"

# ---------------------------------------------------------------------------
# 6. Apply the "Input" cell style to the data rows of both new sheets. Only
#    the cells that actually hold a value get styled (matches the source
#    file, which has no empty-but-styled cells).
# ---------------------------------------------------------------------------
$wsQuestName.Range("A2:C3").Style = "Input"

$wsQuestDesc.Range("A2").Style = "Input"
$wsQuestDesc.Range("B2").Style = "Input"
$wsQuestDesc.Range("E2").Style = "Input"

$wsQuestDesc.Range("A3").Style = "Input"
$wsQuestDesc.Range("B3").Style = "Input"
$wsQuestDesc.Range("C3").Style = "Input"
$wsQuestDesc.Range("E3").Style = "Input"
$wsQuestDesc.Range("F3").Style = "Input"

# ---------------------------------------------------------------------------
# 7. Approximate the bestFit column widths on QuestDescriptionEntities.
# ---------------------------------------------------------------------------
$wsQuestDesc.Columns.Item(1).ColumnWidth = 5.1666666666667
$wsQuestDesc.Columns.Item(2).ColumnWidth = 16
$wsQuestDesc.Columns.Item(3).ColumnWidth = 36
$wsQuestDesc.Columns.Item(4).ColumnWidth = 4.1666666666667
$wsQuestDesc.Columns.Item(5).ColumnWidth = 11.5
$wsQuestDesc.Columns.Item(6).ColumnWidth = 3.6666666666667
$wsQuestDesc.Columns.Item(7).ColumnWidth = 3.6666666666667

# ---------------------------------------------------------------------------
# 8. Restore selections: QuestDescriptionEntities, then QuestNameEntities,
#    then DescriptionEntities, then NameEntities last so NameEntities ends
#    up as the active tab (matches the target file).
# ---------------------------------------------------------------------------
$wsQuestDesc.Activate() | Out-Null
$wsQuestDesc.Range("H14").Select() | Out-Null

$wsQuestName.Activate() | Out-Null
$wsQuestName.Range("A2:XFD3").Select() | Out-Null

$wsDesc.Activate() | Out-Null
$wsDesc.Range("B1:C11").Select() | Out-Null

$wsName.Activate() | Out-Null
$wsName.Range("B15").Select() | Out-Null
